# [Fonds de solidarite] Add 2022-05-30 data
# Update the "nombre_aides" (column C) and "montant_total" (column E) values
# for the rows affected by the new daily data extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 6;   C = 20794;  E = 360648950  },
    @{ Row = 8;   C = 1049;   E = 91299368   },
    @{ Row = 24;  C = 35709;  E = 132410803  },
    @{ Row = 38;  C = 7251;   E = 58559613   },
    @{ Row = 58;  C = 393;    E = 34758363   },
    @{ Row = 70;  C = 15735;  E = 24685528   },
    @{ Row = 92;  C = 409232; E = 1596728283 },
    @{ Row = 93;  C = 209636; E = 1309758854 },
    @{ Row = 94;  C = 94224;  E = 918545109  },
    @{ Row = 96;  C = 17310;  E = 795857580  },
    @{ Row = 104; C = 135256; E = 272260903  },
    @{ Row = 114; C = 3803;   E = 9118147    },
    @{ Row = 120; C = 54;     E = 2534864    },
    @{ Row = 141; C = 80476;  E = 280739870  },
    @{ Row = 144; C = 24417;  E = 201985320  },
    @{ Row = 176; C = 28901;  E = 263226900  },
    @{ Row = 179; C = 634;    E = 54103701   }
)

foreach ($u in $updates) {
    $ws.Range("C" + $u.Row).Value = $u.C
    $ws.Range("E" + $u.Row).Value = $u.E
}
